$wb = $excel.ActiveWorkbook

# 1) Rename worksheets (Product-X -> DescriptiveName-X)
$wb.Worksheets.Item(1).Name = 'NoteBinder-Paper'
$wb.Worksheets.Item(2).Name = 'PrinterPack-Ink'
$wb.Worksheets.Item(3).Name = 'LabelMaker-Screen'
$wb.Worksheets.Item(4).Name = 'BankersBoxSet-Lid'
$wb.Worksheets.Item(5).Name = 'KeyboardTray-Mouse'
$wb.Worksheets.Item(6).Name = 'PencilSharpener-Pencils'

# 2) Update the product-name header cell (C2) on each sheet with the full product name.
#    Order matches the original authoring session (Paper, Ink, Mouse, Screen, Pencils, Lid)
#    so the shared-string table layout mirrors the source workbook byte-for-byte.
$wb.Worksheets.Item(1).Range("C2").Value = 'Five Star Flex Hybrid NoteBinder: Paper'
$wb.Worksheets.Item(2).Range("C2").Value = 'Canon PGI-225 BK/CLI-226 Pack: Ink'
$wb.Worksheets.Item(5).Range("C2").Value = 'Fellowes Designer Suites Premium Keyboard Tray: Mouse'
$wb.Worksheets.Item(3).Range("C2").Value = 'DYMO 1790417 500TS Touchscreen Handheld Label Maker: Screen'
$wb.Worksheets.Item(6).Range("C2").Value = 'X-ACTO SharpX Classic Electric Pencil Sharpener: Pencils'

# 3) Populate the previously-empty summary cells on the Pencils sheet (sheet 6)
$wb.Worksheets.Item(6).Range("C5").Value = 'granted all the pencils were standard #2s, so i am unable to attest to the efficiency of colored pencils or harder/softer leads, but it met and exceeded our needs. it''s very fast and sharpens pencils at a great angle--not too long or short. pencils. )it takes about 20 pencils sharpened to fill the shavings'
$wb.Worksheets.Item(6).Range("D5").Value = 'that occasionally occurs. granted all the pencils were standard #2s, so i am unable to attest to the efficiency of colored pencils or harder/softer leads, but it met and exceeded our needs. it''s very fast and sharpens pencils at a great angle--not too long or short.'
$wb.Worksheets.Item(6).Range("E5").Value = 'granted all the pencils were standard #2s, so i am unable to attest to the efficiency of colored pencils or harder/softer leads, but it met and exceeded our needs. it''s very fast and sharpens pencils at a great angle--not too long or short.'
$wb.Worksheets.Item(6).Range("G5").Value = 'i did find that some pencils took forever to sharpen. it sharpens pencils (and the instructions are very clear that it only sharpens pencils. )it takes about 20 pencils sharpened to fill the shavings chamber. it only does this on standard #2'

$wb.Worksheets.Item(4).Range("C2").Value = 'Bankers Box SmoothMove Moving and Storage Boxes, Small, 10 Pack: Lid'

# 4) Update the selected range on each sheet to match the new authoring state
$wb.Worksheets.Item(1).Range("C2:G2").Select() | Out-Null
$wb.Worksheets.Item(2).Range("C2:G2").Select() | Out-Null
$wb.Worksheets.Item(3).Range("C2:G2").Select() | Out-Null
$wb.Worksheets.Item(4).Range("C2:G2").Select() | Out-Null
$wb.Worksheets.Item(5).Range("C2:G2").Select() | Out-Null
$wb.Worksheets.Item(6).Range("G5").Select() | Out-Null

# 5) Re-activate sheet 1 (tab that was originally selected) and reset the application window geometry
$wb.Worksheets.Item(1).Activate()
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Width = 25600
$excel.ActiveWindow.Height = 16000
